$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap province names and "Casos activos" values between rows 53 and 54
# (Huelva/72 <-> Huesca/0)
$a53 = $ws.Range("A53").Value()
$a54 = $ws.Range("A54").Value()
$c53 = $ws.Range("C53").Value()
$c54 = $ws.Range("C54").Value()

$ws.Range("A53").Value = $a54
$ws.Range("A54").Value = $a53
$ws.Range("C53").Value = $c54
$ws.Range("C54").Value = $c53

# Update the "last updated" timestamp from 07:16 to 07:46
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 07:46"

$wb.Save()
